# Refresh Universalis market-price snapshots + derived leve-profit figures.
# (scheduled runner: pulls latest currentAveragePrice* / LevePrice* / LeveProfit* per leve row)
$wb = $excel.ActiveWorkbook

## Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")

# Row 62: Enchanted Mythrite Ink
$ws.Range("H62").Value = 3616.9
$ws.Range("I62").Value = 2601
$ws.Range("J62").Value = 4632.8
$ws.Range("K62").Value = 2601
$ws.Range("L62").Value = 4632.8
$ws.Range("M62").Value = -1977
$ws.Range("N62").Value = -5880.8

# Row 65: Enchanted Mythrite Ink
$ws.Range("H65").Value = 3616.9
$ws.Range("I65").Value = 2601
$ws.Range("J65").Value = 4632.8
$ws.Range("K65").Value = 13005
$ws.Range("L65").Value = 23164
$ws.Range("M65").Value = -9885
$ws.Range("N65").Value = -29404

# Row 70: Holy Water
$ws.Range("H70").Value = 3201.7727
$ws.Range("I70").Value = 1326.6666
$ws.Range("J70").Value = 3497.842
$ws.Range("K70").Value = 3979.9998
$ws.Range("L70").Value = 10493.526
$ws.Range("M70").Value = -3709.9998
$ws.Range("N70").Value = -11033.526

# Row 73: Holy Water
$ws.Range("H73").Value = 3201.7727
$ws.Range("I73").Value = 1326.6666
$ws.Range("J73").Value = 3497.842
$ws.Range("K73").Value = 3979.9998
$ws.Range("L73").Value = 10493.526
$ws.Range("M73").Value = -3043.9998
$ws.Range("N73").Value = -12365.526

# Row 135: Grade 1 Gemsap of Intelligence
$ws.Range("H135").Value = 686.08
$ws.Range("I135").Value = 697.8182
$ws.Range("K135").Value = 6280.3638
$ws.Range("M135").Value = -3745.3638

# Row 137: Magnesia Whetstone
$ws.Range("H137").Value = 3461.5667
$ws.Range("I137").Value = 4079.1765
$ws.Range("K137").Value = 12237.5295
$ws.Range("M137").Value = -9687.529500000001

# Row 138: Cunning Craftsman's Tisane
$ws.Range("H138").Value = 3763.7793
$ws.Range("I138").Value = 1681.5682
$ws.Range("J138").Value = 7581.1665
$ws.Range("K138").Value = 5044.7046
$ws.Range("L138").Value = 22743.4995
$ws.Range("M138").Value = 95.29539999999997
$ws.Range("N138").Value = -33023.49950000001

## Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")

# Row 2: Bronze Ingot
$ws.Range("H2").Value = 4903501
$ws.Range("I2").Value = 10417435
$ws.Range("J2").Value = 2226.2222
$ws.Range("K2").Value = 10417435
$ws.Range("L2").Value = 2226.2222
$ws.Range("M2").Value = -10417322
$ws.Range("N2").Value = -2452.2222

# Row 116: Titanbronze Ingot
$ws.Range("H116").Value = 4903501
$ws.Range("I116").Value = 10417435
$ws.Range("J116").Value = 2226.2222
$ws.Range("K116").Value = 10417435
$ws.Range("L116").Value = 2226.2222
$ws.Range("M116").Value = -10415141
$ws.Range("N116").Value = -6814.2222

# Row 122: High Durium Nugget
$ws.Range("H122").Value = 2181.5833
$ws.Range("I122").Value = 1416.5
$ws.Range("J122").Value = 2946.6667
$ws.Range("K122").Value = 4249.5
$ws.Range("L122").Value = 8840.000100000001
$ws.Range("M122").Value = -1799.5
$ws.Range("N122").Value = -13740.0001

# Row 132: Mountain Chromite Ingot
$ws.Range("H132").Value = 14288348
$ws.Range("I132").Value = 18870432
$ws.Range("J132").Value = 3028.2354
$ws.Range("K132").Value = 56611296
$ws.Range("L132").Value = 9084.706200000001
$ws.Range("M132").Value = -56608766
$ws.Range("N132").Value = -14144.7062

## Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")

# Row 3: Bronze Ingot
$ws.Range("H3").Value = 4903501
$ws.Range("I3").Value = 10417435
$ws.Range("J3").Value = 2226.2222
$ws.Range("K3").Value = 10417435
$ws.Range("L3").Value = 2226.2222
$ws.Range("M3").Value = -10417321
$ws.Range("N3").Value = -2454.2222

# Row 20: Iron Ingot
$ws.Range("H20").Value = 1485.7826
$ws.Range("I20").Value = 718.6667
$ws.Range("J20").Value = 2322.6365
$ws.Range("K20").Value = 718.6667
$ws.Range("L20").Value = 2322.6365
$ws.Range("M20").Value = -471.6667
$ws.Range("N20").Value = -2816.6365

# Row 86: Adamantite Nugget
$ws.Range("H86").Value = 2830
$ws.Range("I86").Value = 2254.1667
$ws.Range("J86").Value = 5133.3335
$ws.Range("K86").Value = 2254.1667
$ws.Range("L86").Value = 5133.3335
$ws.Range("M86").Value = -1131.1667
$ws.Range("N86").Value = -7379.3335

# Row 89: Adamantite Nugget
$ws.Range("H89").Value = 2830
$ws.Range("I89").Value = 2254.1667
$ws.Range("J89").Value = 5133.3335
$ws.Range("K89").Value = 11270.8335
$ws.Range("L89").Value = 25666.6675
$ws.Range("M89").Value = -5654.833500000001
$ws.Range("N89").Value = -36898.6675

# Row 99: Oroshigane Ingot
$ws.Range("H99").Value = 1492.0869
$ws.Range("I99").Value = 1048.3158
$ws.Range("J99").Value = 3600
$ws.Range("K99").Value = 1048.3158
$ws.Range("L99").Value = 3600
$ws.Range("M99").Value = 449.6841999999999
$ws.Range("N99").Value = -6596

# Row 134: Ruthenium Ingot
$ws.Range("H134").Value = 2504.8572
$ws.Range("I134").Value = 2188.756
$ws.Range("J134").Value = 4124.875
$ws.Range("K134").Value = 6566.268
$ws.Range("L134").Value = 12374.625
$ws.Range("M134").Value = -4031.268
$ws.Range("N134").Value = -17444.625

## Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")

# Row 31: Walnut Lumber
$ws.Range("H31").Value = 2351.418
$ws.Range("I31").Value = 1488.9736
$ws.Range("J31").Value = 3481.5173
$ws.Range("K31").Value = 1488.9736
$ws.Range("L31").Value = 3481.5173
$ws.Range("M31").Value = -1193.9736
$ws.Range("N31").Value = -4071.5173

# Row 34: Walnut Lumber
$ws.Range("H34").Value = 2351.418
$ws.Range("I34").Value = 1488.9736
$ws.Range("J34").Value = 3481.5173
$ws.Range("K34").Value = 1488.9736
$ws.Range("L34").Value = 3481.5173
$ws.Range("M34").Value = -1286.9736
$ws.Range("N34").Value = -3885.5173

# Row 132: Ginseng Lumber
$ws.Range("H132").Value = 2199.8157
$ws.Range("I132").Value = 1313.4333
$ws.Range("J132").Value = 5523.75
$ws.Range("K132").Value = 3940.2999
$ws.Range("L132").Value = 16571.25
$ws.Range("M132").Value = -1410.2999
$ws.Range("N132").Value = -21631.25

# Row 134: Ceiba Lumber
$ws.Range("H134").Value = 1560.3433
$ws.Range("I134").Value = 1207.849
$ws.Range("K134").Value = 3623.547
$ws.Range("M134").Value = -1088.547

## Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")

# Row 12: Kukuru Butter
$ws.Range("H12").Value = 215.26666
$ws.Range("I12").Value = 29.5
$ws.Range("J12").Value = 243.84616
$ws.Range("K12").Value = 88.5
$ws.Range("L12").Value = 731.5384799999999
$ws.Range("M12").Value = 84.5
$ws.Range("N12").Value = -1077.53848

# Row 68: Fermented Butter
$ws.Range("H68").Value = 2200.3867
$ws.Range("I68").Value = 735.0333000000001
$ws.Range("J68").Value = 3177.2888
$ws.Range("K68").Value = 2205.0999
$ws.Range("L68").Value = 9531.866399999999
$ws.Range("M68").Value = -1394.0999
$ws.Range("N68").Value = -11153.8664

# Row 71: Fermented Butter
$ws.Range("H71").Value = 2200.3867
$ws.Range("I71").Value = 735.0333000000001
$ws.Range("J71").Value = 3177.2888
$ws.Range("K71").Value = 6615.2997
$ws.Range("L71").Value = 28595.5992
$ws.Range("M71").Value = -2559.2997
$ws.Range("N71").Value = -36707.5992

# Row 131: Tsai tou Vounou
$ws.Range("H131").Value = 1219.5217
$ws.Range("I131").Value = 1222.1428
$ws.Range("J131").Value = 1218.375
$ws.Range("K131").Value = 3666.4284
$ws.Range("L131").Value = 3655.125
$ws.Range("M131").Value = 1373.5716
$ws.Range("N131").Value = -13735.125

## Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")

# Row 97: Koppranickel Ingot
$ws.Range("H97").Value = 1608.381
$ws.Range("I97").Value = 1222.1177
$ws.Range("J97").Value = 3250
$ws.Range("K97").Value = 1222.1177
$ws.Range("L97").Value = 3250
$ws.Range("M97").Value = -726.1177
$ws.Range("N97").Value = -4242

# Row 122: Ametrine
$ws.Range("H122").Value = 8250
$ws.Range("I122").Value = 15000
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 45000
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -42550
$ws.Range("N122").Value = -22900

# Row 132: Lar Ingot
$ws.Range("H132").Value = 4375
$ws.Range("I132").Value = 3375
$ws.Range("K132").Value = 10125
$ws.Range("M132").Value = -7595

## Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")

# Row 55: Peiste Leather
$ws.Range("H55").Value = 979.5625
$ws.Range("I55").Value = 255.8
$ws.Range("J55").Value = 1308.5454
$ws.Range("K55").Value = 255.8
$ws.Range("L55").Value = 1308.5454
$ws.Range("M55").Value = -82.80000000000001
$ws.Range("N55").Value = -1654.5454

# Row 132: Silver Lobo Leather
$ws.Range("H132").Value = 1950.1936
$ws.Range("J132").Value = 2863.25
$ws.Range("L132").Value = 8589.75
$ws.Range("N132").Value = -13649.75

# Row 138: Gomphotherium Boots of Striking
$ws.Range("H138").Value = 40429
$ws.Range("J138").Value = 40429
$ws.Range("L138").Value = 40429
$ws.Range("N138").Value = -50709

## Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")

# Row 132: Snow Cotton Cloth
$ws.Range("H132").Value = 7942.1665
$ws.Range("I132").Value = 3300.5386
$ws.Range("J132").Value = 13427.728
$ws.Range("K132").Value = 9901.6158
$ws.Range("L132").Value = 40283.18399999999
$ws.Range("M132").Value = -7371.6158
$ws.Range("N132").Value = -45343.18399999999
